$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted style from an existing "Date" cell (column A)
# down to the new row so the new date cell keeps the same number format
# (style index) instead of Excel minting a brand-new style entry.
$ws.Range("A45").Copy($ws.Range("A49"))

# New log entry (row 49): 11/22/2019, 10:30AM - 12:00PM, 90 min delta,
# with a note about the new generic messageBox confirmation/error signal.
$ws.Range("A49").Value = 43791
$ws.Range("B49").Value = "10:30AM"
$ws.Range("C49").Value = "12:00PM"
$ws.Range("E49").Value = 90
$ws.Range("F49").Value = "Setting up generic messageBox to show confirmation or error msg to user"

# Update the active selection to reflect where the user ended up after
# adding the new row.
$ws.Activate()
$ws.Range("F50").Select()
